$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Price (column D) text values
$dValues = @{
    2 = "24.430.84"
    3 = "1.655.95"
    4 = "1.001"
    5 = "311.71"
    6 = "1.001"
    7 = "0.3918"
    8 = "0.3905"
    9 = "1.002"
    10 = "1.382"
    11 = "50.06"
    12 = "0.08562"
    13 = "24.96"
    15 = "0.00001305"
    16 = "7.612"
    17 = "1.658.82"
    18 = "93.23"
    19 = "0.06953"
    20 = "20.97"
    21 = "7.006"
    22 = "1.002"
    23 = "13.81"
    24 = "24.421.27"
    25 = "2.337"
    26 = "2.790"
    27 = "22.70"
    28 = "158.86"
    29 = "5.730"
    30 = "145.36"
    31 = "8.190"
    32 = "2.518"
    33 = "1.836.51"
    34 = "0.03017"
    36 = "1.002"
    38 = "0.2765"
    39 = "0.09468"
    40 = "1.491"
    41 = "10.20"
    42 = "0.7784"
    43 = "13.35"
    44 = "16.38"
    45 = "2.552"
    46 = "0.7017"
    47 = "4.146"
    48 = "0.08585"
    49 = "1.001"
    50 = "1.304"
    51 = "136.44"
}

# Map of row -> new Volume(1h) (column E) text values
$eValues = @{
    2 = "  -1.37%  "
    3 = "  -2.52%  "
    4 = "  -0.07%  "
    6 = "  +0.00%  "
    7 = "  -1.64%  "
    9 = "  -0.08%  "
    10 = "  -5.76%  "
    11 = "  -6.53%  "
    12 = "  -2.67%  "
    13 = "  -5.29%  "
    14 = "  -4.41%  "
    15 = "  -2.69%  "
    16 = "  -4.58%  "
    17 = "  -4.41%  "
    18 = "  -2.37%  "
    19 = "  -3.12%  "
    20 = "  +0.56%  "
    21 = "  -4.48%  "
    22 = "  -0.02%  "
    23 = "  -3.98%  "
    24 = "  -1.29%  "
    25 = "  -1.89%  "
    26 = "  -3.97%  "
    27 = "  -1.62%  "
    28 = "  -2.00%  "
    29 = "  -6.34%  "
    30 = "  +0.65%  "
    31 = "  -2.00%  "
    32 = "  +11.48%  "
    33 = "  -2.17%  "
    34 = "  -5.03%  "
    35 = "  -5.78%  "
    36 = "  -2.46%  "
    37 = "  -6.47%  "
    38 = "  -2.88%  "
    39 = "  +0.50%  "
    40 = "  +0.94%  "
    41 = "  -4.73%  "
    42 = "  -6.32%  "
    43 = "  -5.56%  "
    44 = "  -6.92%  "
    45 = "  -5.64%  "
    46 = "  -5.45%  "
    47 = "  -1.84%  "
    48 = "  +2.46%  "
    49 = "  +0.03%  "
    50 = "  -5.34%  "
    51 = "  -2.14%  "
}

foreach ($row in $dValues.Keys) {
    $cell = $ws.Range("D$row")
    # Force the cell to Text format first so Excel does not auto-convert
    # numeric-looking strings (e.g. "1.001", "311.71") into actual numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $dValues[$row]
    # Restore the default "Normal" style so no stray style index is left on the cell.
    $cell.Style = "Normal"
}

foreach ($row in $eValues.Keys) {
    $ws.Range("E$row").Value = $eValues[$row]
}
